# training_info.xlsx update
#  - "increased penalty for drinks": the training-instance count for the
#    "alter" / "that was too <flavor>" row (F7) goes up from 44 to 50.
#  - "fixed confirmation": leave the active selection resting on the cell
#    that was actually edited (F7) instead of the stray H7 selection that
#    was saved previously.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the training-instance count in F7 (row for "alter" / flavor=N/A?).
$ws.Cells.Item(7, 6).Value = 50

# Move/confirm the selection onto the cell that was just edited.
$ws.Range("F7").Select() | Out-Null
